# Auto-generated edit script applying TPM value updates to L1cam-Egfr sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = [double]"6.603177"
$ws.Range("H2").Value = [double]"19.809531"
$ws.Range("I2").Value = [double]"0.5135477412645301"
$ws.Range("J2").Value = [double]"0.5135477412645302"
$ws.Range("M2").Value = [double]"1.370876333333333"
$ws.Range("N2").Value = [double]"4.112629"
$ws.Range("O2").Value = [double]"0.01103063309339269"
$ws.Range("P2").Value = [double]"0.01103063309339269"
$ws.Range("Q2").Value = [double]"9.052139074111"
$ws.Range("R2").Value = [double]"81.469251666999"
$ws.Range("S2").Value = [double]"0.005664756709829594"
$ws.Range("T2").Value = [double]"0.005664756709829594"
$ws.Range("G3").Value = [double]"6.603177"
$ws.Range("H3").Value = [double]"19.809531"
$ws.Range("I3").Value = [double]"0.5135477412645301"
$ws.Range("J3").Value = [double]"0.5135477412645302"
$ws.Range("O3").Value = [double]"0.7476219244149905"
$ws.Range("P3").Value = [double]"0.7476219244149904"
$ws.Range("Q3").Value = [double]"613.5257675022069"
$ws.Range("R3").Value = [double]"5521.731907519863"
$ws.Range("S3").Value = [double]"0.3839395506031596"
$ws.Range("T3").Value = [double]"0.3839395506031597"
$ws.Range("G4").Value = [double]"6.603177"
$ws.Range("H4").Value = [double]"19.809531"
$ws.Range("I4").Value = [double]"0.5135477412645301"
$ws.Range("J4").Value = [double]"0.5135477412645302"
$ws.Range("M4").Value = [double]"29.718484"
$ws.Range("N4").Value = [double]"89.155452"
$ws.Range("O4").Value = [double]"0.2391271080585153"
$ws.Range("P4").Value = [double]"0.2391271080585153"
$ws.Range("Q4").Value = [double]"196.236410023668"
$ws.Range("R4").Value = [double]"1766.127690213012"
$ws.Range("S4").Value = [double]"0.1228031862185698"
$ws.Range("T4").Value = [double]"0.1228031862185698"
$ws.Range("G5").Value = [double]"6.603177"
$ws.Range("H5").Value = [double]"19.809531"
$ws.Range("I5").Value = [double]"0.5135477412645301"
$ws.Range("J5").Value = [double]"0.5135477412645302"
$ws.Range("M5").Value = [double]"0.275941"
$ws.Range("N5").Value = [double]"0.827823"
$ws.Range("O5").Value = [double]"0.002220334433101459"
$ws.Range("P5").Value = [double]"0.002220334433101458"
$ws.Range("Q5").Value = [double]"1.822087264557"
$ws.Range("R5").Value = [double]"16.398785381013"
$ws.Range("S5").Value = [double]"0.001140247732971115"
$ws.Range("T5").Value = [double]"0.001140247732971115"
$ws.Range("I6").Value = [double]"0.02944398858046029"
$ws.Range("J6").Value = [double]"0.0294439885804603"
$ws.Range("M6").Value = [double]"1.370876333333333"
$ws.Range("N6").Value = [double]"4.112629"
$ws.Range("O6").Value = [double]"0.01103063309339269"
$ws.Range("P6").Value = [double]"0.01103063309339269"
$ws.Range("Q6").Value = [double]"0.5189996140778889"
$ws.Range("R6").Value = [double]"4.670996526701"
$ws.Range("S6").Value = [double]"0.0003247858348371019"
$ws.Range("T6").Value = [double]"0.0003247858348371019"
$ws.Range("I7").Value = [double]"0.02944398858046029"
$ws.Range("J7").Value = [double]"0.0294439885804603"
$ws.Range("O7").Value = [double]"0.7476219244149905"
$ws.Range("P7").Value = [double]"0.7476219244149904"
$ws.Range("S7").Value = [double]"0.02201297140497673"
$ws.Range("T7").Value = [double]"0.02201297140497673"
$ws.Range("I8").Value = [double]"0.02944398858046029"
$ws.Range("J8").Value = [double]"0.0294439885804603"
$ws.Range("M8").Value = [double]"29.718484"
$ws.Range("N8").Value = [double]"89.155452"
$ws.Range("O8").Value = [double]"0.2391271080585153"
$ws.Range("P8").Value = [double]"0.2391271080585153"
$ws.Range("Q8").Value = [double]"11.25111095139867"
$ws.Range("R8").Value = [double]"101.259998562588"
$ws.Range("S8").Value = [double]"0.00704085583895342"
$ws.Range("T8").Value = [double]"0.007040855838953419"
$ws.Range("I9").Value = [double]"0.02944398858046029"
$ws.Range("J9").Value = [double]"0.0294439885804603"
$ws.Range("M9").Value = [double]"0.275941"
$ws.Range("N9").Value = [double]"0.827823"
$ws.Range("O9").Value = [double]"0.002220334433101459"
$ws.Range("P9").Value = [double]"0.002220334433101458"
$ws.Range("Q9").Value = [double]"0.1044684112096667"
$ws.Range("R9").Value = [double]"0.940215700887"
$ws.Range("S9").Value = [double]"6.537550169304213E-05"
$ws.Range("T9").Value = [double]"6.537550169304212E-05"
$ws.Range("G10").Value = [double]"3.441487333333333"
$ws.Range("H10").Value = [double]"10.324462"
$ws.Range("I10").Value = [double]"0.2676541983690312"
$ws.Range("J10").Value = [double]"0.2676541983690313"
$ws.Range("M10").Value = [double]"1.370876333333333"
$ws.Range("N10").Value = [double]"4.112629"
$ws.Range("O10").Value = [double]"0.01103063309339269"
$ws.Range("P10").Value = [double]"0.01103063309339269"
$ws.Range("Q10").Value = [double]"4.717853536733111"
$ws.Range("R10").Value = [double]"42.46068183059801"
$ws.Range("S10").Value = [double]"0.002952395258114928"
$ws.Range("T10").Value = [double]"0.002952395258114929"
$ws.Range("G11").Value = [double]"3.441487333333333"
$ws.Range("H11").Value = [double]"10.324462"
$ws.Range("I11").Value = [double]"0.2676541983690312"
$ws.Range("J11").Value = [double]"0.2676541983690313"
$ws.Range("O11").Value = [double]"0.7476219244149905"
$ws.Range("P11").Value = [double]"0.7476219244149904"
$ws.Range("Q11").Value = [double]"319.7614053859918"
$ws.Range("R11").Value = [double]"2877.852648473926"
$ws.Range("S11").Value = [double]"0.2001041468624067"
$ws.Range("T11").Value = [double]"0.2001041468624067"
$ws.Range("G12").Value = [double]"3.441487333333333"
$ws.Range("H12").Value = [double]"10.324462"
$ws.Range("I12").Value = [double]"0.2676541983690312"
$ws.Range("J12").Value = [double]"0.2676541983690313"
$ws.Range("M12").Value = [double]"29.718484"
$ws.Range("N12").Value = [double]"89.155452"
$ws.Range("O12").Value = [double]"0.2391271080585153"
$ws.Range("P12").Value = [double]"0.2391271080585153"
$ws.Range("Q12").Value = [double]"102.2757862518693"
$ws.Range("R12").Value = [double]"920.482076266824"
$ws.Range("S12").Value = [double]"0.06400337441570662"
$ws.Range("T12").Value = [double]"0.06400337441570662"
$ws.Range("G13").Value = [double]"3.441487333333333"
$ws.Range("H13").Value = [double]"10.324462"
$ws.Range("I13").Value = [double]"0.2676541983690312"
$ws.Range("J13").Value = [double]"0.2676541983690313"
$ws.Range("M13").Value = [double]"0.275941"
$ws.Range("N13").Value = [double]"0.827823"
$ws.Range("O13").Value = [double]"0.002220334433101459"
$ws.Range("P13").Value = [double]"0.002220334433101458"
$ws.Range("Q13").Value = [double]"0.9496474562473333"
$ws.Range("R13").Value = [double]"8.546827106226001"
$ws.Range("S13").Value = [double]"0.0005942818328029283"
$ws.Range("T13").Value = [double]"0.0005942818328029283"
$ws.Range("G14").Value = [double]"2.434707333333333"
$ws.Range("H14").Value = [double]"7.304122"
$ws.Range("I14").Value = [double]"0.1893540717859783"
$ws.Range("J14").Value = [double]"0.1893540717859783"
$ws.Range("M14").Value = [double]"1.370876333333333"
$ws.Range("N14").Value = [double]"4.112629"
$ws.Range("O14").Value = [double]"0.01103063309339269"
$ws.Range("P14").Value = [double]"0.01103063309339269"
$ws.Range("Q14").Value = [double]"3.337682661859778"
$ws.Range("R14").Value = [double]"30.039143956738"
$ws.Range("S14").Value = [double]"0.002088695290611068"
$ws.Range("T14").Value = [double]"0.002088695290611068"
$ws.Range("G15").Value = [double]"2.434707333333333"
$ws.Range("H15").Value = [double]"7.304122"
$ws.Range("I15").Value = [double]"0.1893540717859783"
$ws.Range("J15").Value = [double]"0.1893540717859783"
$ws.Range("O15").Value = [double]"0.7476219244149905"
$ws.Range("P15").Value = [double]"0.7476219244149904"
$ws.Range("Q15").Value = [double]"226.2177260016784"
$ws.Range("R15").Value = [double]"2035.959534015106"
$ws.Range("S15").Value = [double]"0.1415652555444473"
$ws.Range("T15").Value = [double]"0.1415652555444473"
$ws.Range("G16").Value = [double]"2.434707333333333"
$ws.Range("H16").Value = [double]"7.304122"
$ws.Range("I16").Value = [double]"0.1893540717859783"
$ws.Range("J16").Value = [double]"0.1893540717859783"
$ws.Range("M16").Value = [double]"29.718484"
$ws.Range("N16").Value = [double]"89.155452"
$ws.Range("O16").Value = [double]"0.2391271080585153"
$ws.Range("P16").Value = [double]"0.2391271080585153"
$ws.Range("Q16").Value = [double]"72.35581093034934"
$ws.Range("R16").Value = [double]"651.202298373144"
$ws.Range("S16").Value = [double]"0.0452796915852855"
$ws.Range("T16").Value = [double]"0.0452796915852855"
$ws.Range("G17").Value = [double]"2.434707333333333"
$ws.Range("H17").Value = [double]"7.304122"
$ws.Range("I17").Value = [double]"0.1893540717859783"
$ws.Range("J17").Value = [double]"0.1893540717859783"
$ws.Range("M17").Value = [double]"0.275941"
$ws.Range("N17").Value = [double]"0.827823"
$ws.Range("O17").Value = [double]"0.002220334433101459"
$ws.Range("P17").Value = [double]"0.002220334433101458"
$ws.Range("Q17").Value = [double]"0.6718355762673334"
$ws.Range("R17").Value = [double]"6.046520186406"
$ws.Range("S17").Value = [double]"0.000420429365634373"
$ws.Range("T17").Value = [double]"0.000420429365634373"
